$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.04628467895650446, 1.329827146587384, 9.907280297060181, 3.147583247042114, 3.177948460597182, 52)
    3  = @(0.03072812890343799, 1.487177437779083, 9.489403539255518, 3.080487548953172, 3.110985134611961, 51)
    4  = @(0.05670300065403502, 1.350516681017773, 6.724916598423166, 2.593244415480956, 2.61894615200115, 50)
    5  = @(0.06807619608839131, 1.406782146406506, 8.32307609620424, 2.884974193334187, 2.914059470405175, 49)
    6  = @(0.0465829568975624, 1.552000945598137, 9.764321243549908, 3.124791392005218, 3.15750797590164, 48)
    7  = @(0.007348863100462911, 1.800235720797896, 10.52634790614585, 3.244433372123066, 3.290447561450098, 36)
    8  = @(0.04734274986515632, 1.848121295394991, 10.96964521407898, 3.312045472827778, 3.360055740670333, 35)
    9  = @(0.08958096046880724, 2.597962177260555, 19.61373083157788, 4.428739192092698, 4.55620268551751, 18)
    10 = @(-0.6075801358359439, 3.106648192679968, 23.3436798468976, 4.831529762600827, 5.027124311967897, 11)
    11 = @(-1.322332775507276, 2.379927842974612, 17.91130170936413, 4.232174583989196, 4.494821706723103, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
